# Update main GSC export data:
#  - Drop the oldest day (2025-11-20) from the "Chart" sheet, shifting all
#    subsequent days up by one row.
#  - The two newest days that now lack validated data (2025-11-21 and
#    2025-11-22, plus the following day) show blank "No video indexed" /
#    "Video indexed" counts instead of numbers.
#  - The last row (now 2026-02-16) no longer carries a Reason/Validation
#    text value in column D; it is a plain numeric 0 like the other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the first data row (2025-11-20); everything below shifts up.
$ws.Rows.Item(2).Delete()

# The first three remaining data rows (2025-11-21, 2025-11-22, 2025-11-23)
# no longer have numeric "No video indexed" / "Video indexed" counts.
$ws.Cells.Item(2, 2).Value = ""
$ws.Cells.Item(2, 3).Value = ""
$ws.Cells.Item(3, 2).Value = ""
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(4, 2).Value = ""
$ws.Cells.Item(4, 3).Value = ""

# The new last row (2026-02-16) now has a numeric 0 instead of blank text.
$ws.Cells.Item(89, 4).Value = 0
